$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 62501876
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H115").Value = 2371.7273
$ws.Range("I115").Value = 2254.4443
$ws.Range("J115").Value = 2899.5
$ws.Range("K115").Value = 6763.3329
$ws.Range("L115").Value = 8698.5
$ws.Range("M115").Value = -5196.3329
$ws.Range("N115").Value = -11832.5
$ws.Range("H132").Value = 24733760
$ws.Range("I132").Value = 27407100
$ws.Range("J132").Value = 5375
$ws.Range("K132").Value = 82221300
$ws.Range("L132").Value = 16125
$ws.Range("M132").Value = -82218770
$ws.Range("N132").Value = -21185
$ws.Range("H137").Value = 3770.3726
$ws.Range("I137").Value = 3689.25
$ws.Range("J137").Value = 3907
$ws.Range("K137").Value = 11067.75
$ws.Range("L137").Value = 11721
$ws.Range("M137").Value = -8517.75
$ws.Range("N137").Value = -16821
$ws.Range("H138").Value = 3997.9255
$ws.Range("I138").Value = 2687.889
$ws.Range("J138").Value = 4136.6353
$ws.Range("K138").Value = 8063.667
$ws.Range("L138").Value = 12409.9059
$ws.Range("M138").Value = -2923.667
$ws.Range("N138").Value = -22689.9059
$ws.Range("H141").Value = 1949.0333
$ws.Range("I141").Value = 1519.64
$ws.Range("J141").Value = 4096
$ws.Range("K141").Value = 4558.92
$ws.Range("L141").Value = 12288
$ws.Range("M141").Value = 621.0799999999999
$ws.Range("N141").Value = -22648

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13947.202
$ws.Range("I32").Value = 9601
$ws.Range("J32").Value = 20357.85
$ws.Range("K32").Value = 9601
$ws.Range("L32").Value = 20357.85
$ws.Range("M32").Value = -9314
$ws.Range("N32").Value = -20931.85
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 10500
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 10500
$ws.Range("M43").Value = -4187
$ws.Range("N43").Value = -11126
$ws.Range("H61").Value = 3712.4443
$ws.Range("I61").Value = 3058.8572
$ws.Range("K61").Value = 3058.8572
$ws.Range("M61").Value = -2846.8572
$ws.Range("H74").Value = 4158.4546
$ws.Range("I74").Value = 4635.8696
$ws.Range("J74").Value = 3060.4
$ws.Range("K74").Value = 4635.8696
$ws.Range("L74").Value = 3060.4
$ws.Range("M74").Value = -3761.8696
$ws.Range("N74").Value = -4808.4
$ws.Range("H77").Value = 4158.4546
$ws.Range("I77").Value = 4635.8696
$ws.Range("J77").Value = 3060.4
$ws.Range("K77").Value = 23179.348
$ws.Range("L77").Value = 15302
$ws.Range("M77").Value = -18811.348
$ws.Range("N77").Value = -24038
$ws.Range("H97").Value = 1066.3334
$ws.Range("I97").Value = 910.7222
$ws.Range("K97").Value = 910.7222
$ws.Range("M97").Value = -414.7222
$ws.Range("H110").Value = 637.6
$ws.Range("I110").Value = 661.8333
$ws.Range("J110").Value = 540.6667
$ws.Range("K110").Value = 661.8333
$ws.Range("L110").Value = 540.6667
$ws.Range("M110").Value = 1383.1667
$ws.Range("N110").Value = -4630.6667
$ws.Range("H122").Value = 2837.1765
$ws.Range("I122").Value = 1804
$ws.Range("J122").Value = 6195
$ws.Range("K122").Value = 5412
$ws.Range("L122").Value = 18585
$ws.Range("M122").Value = -2962
$ws.Range("N122").Value = -23485
$ws.Range("H136").Value = 3712.4443
$ws.Range("I136").Value = 3058.8572
$ws.Range("K136").Value = 9176.571599999999
$ws.Range("M136").Value = -6626.571599999999
$ws.Range("H137").Value = 44390
$ws.Range("J137").Value = 44390
$ws.Range("L137").Value = 44390
$ws.Range("N137").Value = -54590

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10478.579
$ws.Range("I20").Value = 4798.25
$ws.Range("J20").Value = 14609.728
$ws.Range("K20").Value = 4798.25
$ws.Range("L20").Value = 14609.728
$ws.Range("M20").Value = -4551.25
$ws.Range("N20").Value = -15103.728
$ws.Range("H94").Value = 17858656
$ws.Range("I94").Value = 23810852
$ws.Range("J94").Value = 2070.1428
$ws.Range("K94").Value = 23810852
$ws.Range("L94").Value = 2070.1428
$ws.Range("M94").Value = -23810401
$ws.Range("N94").Value = -2972.1428
$ws.Range("H134").Value = 3762.7942
$ws.Range("I134").Value = 2034.1111
$ws.Range("J134").Value = 10430.571
$ws.Range("K134").Value = 6102.3333
$ws.Range("L134").Value = 31291.713
$ws.Range("M134").Value = -3567.3333
$ws.Range("N134").Value = -36361.713

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3962.4897
$ws.Range("I31").Value = 1488.7391
$ws.Range("J31").Value = 6150.8076
$ws.Range("K31").Value = 1488.7391
$ws.Range("L31").Value = 6150.8076
$ws.Range("M31").Value = -1193.7391
$ws.Range("N31").Value = -6740.8076
$ws.Range("H34").Value = 3962.4897
$ws.Range("I34").Value = 1488.7391
$ws.Range("J34").Value = 6150.8076
$ws.Range("K34").Value = 1488.7391
$ws.Range("L34").Value = 6150.8076
$ws.Range("M34").Value = -1286.7391
$ws.Range("N34").Value = -6554.8076
$ws.Range("H122").Value = 4013.7856
$ws.Range("I122").Value = 2456.1428
$ws.Range("K122").Value = 7368.428400000001
$ws.Range("M122").Value = -4918.428400000001
$ws.Range("H132").Value = 2929.1667
$ws.Range("I132").Value = 2330.8057
$ws.Range("J132").Value = 4724.25
$ws.Range("K132").Value = 6992.4171
$ws.Range("L132").Value = 14172.75
$ws.Range("M132").Value = -4462.4171
$ws.Range("N132").Value = -19232.75
$ws.Range("H134").Value = 5456.037
$ws.Range("I134").Value = 5875.727
$ws.Range("J134").Value = 3609.4
$ws.Range("K134").Value = 17627.181
$ws.Range("L134").Value = 10828.2
$ws.Range("M134").Value = -15092.181
$ws.Range("N134").Value = -15898.2

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1994.96
$ws.Range("I5").Value = 618.26666
$ws.Range("J5").Value = 4060
$ws.Range("K5").Value = 1854.79998
$ws.Range("L5").Value = 12180
$ws.Range("M5").Value = -1742.79998
$ws.Range("N5").Value = -12404
$ws.Range("H132").Value = 2242.55
$ws.Range("I132").Value = 1005.35297
$ws.Range("K132").Value = 9048.176730000001
$ws.Range("M132").Value = -6518.176730000001
$ws.Range("H135").Value = 1994.96
$ws.Range("I135").Value = 618.26666
$ws.Range("J135").Value = 4060
$ws.Range("K135").Value = 5564.39994
$ws.Range("L135").Value = 36540
$ws.Range("M135").Value = -3029.39994
$ws.Range("N135").Value = -41610

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5992.5293
$ws.Range("I70").Value = 5470.0454
$ws.Range("J70").Value = 6950.4165
$ws.Range("K70").Value = 5470.0454
$ws.Range("L70").Value = 6950.4165
$ws.Range("M70").Value = -5200.0454
$ws.Range("N70").Value = -7490.4165
$ws.Range("H73").Value = 5992.5293
$ws.Range("I73").Value = 5470.0454
$ws.Range("J73").Value = 6950.4165
$ws.Range("K73").Value = 5470.0454
$ws.Range("L73").Value = 6950.4165
$ws.Range("M73").Value = -4534.0454
$ws.Range("N73").Value = -8822.416499999999
$ws.Range("H97").Value = 1366.875
$ws.Range("I97").Value = 1302.8182
$ws.Range("J97").Value = 1507.8
$ws.Range("K97").Value = 1302.8182
$ws.Range("L97").Value = 1507.8
$ws.Range("M97").Value = -806.8181999999999
$ws.Range("N97").Value = -2499.8
$ws.Range("H102").Value = 2699.38
$ws.Range("I102").Value = 2333.8372
$ws.Range("J102").Value = 4944.857
$ws.Range("K102").Value = 2333.8372
$ws.Range("L102").Value = 4944.857
$ws.Range("M102").Value = -711.8371999999999
$ws.Range("N102").Value = -8188.857
$ws.Range("H132").Value = 2336.468
$ws.Range("I132").Value = 1123.7727
$ws.Range("J132").Value = 3403.64
$ws.Range("K132").Value = 3371.3181
$ws.Range("L132").Value = 10210.92
$ws.Range("M132").Value = -841.3181
$ws.Range("N132").Value = -15270.92
$ws.Range("H137").Value = 61061.145
$ws.Range("J137").Value = 61061.145
$ws.Range("L137").Value = 61061.145
$ws.Range("N137").Value = -71261.14499999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2542.8572
$ws.Range("I22").Value = 2090.818
$ws.Range("J22").Value = 3040.1
$ws.Range("K22").Value = 2090.818
$ws.Range("L22").Value = 3040.1
$ws.Range("M22").Value = -1795.818
$ws.Range("N22").Value = -3630.1
$ws.Range("H27").Value = 2542.8572
$ws.Range("I27").Value = 2090.818
$ws.Range("J27").Value = 3040.1
$ws.Range("K27").Value = 2090.818
$ws.Range("L27").Value = 3040.1
$ws.Range("M27").Value = -1983.818
$ws.Range("N27").Value = -3254.1
$ws.Range("H100").Value = 2038.7
$ws.Range("I100").Value = 1629.8334
$ws.Range("K100").Value = 1629.8334
$ws.Range("M100").Value = -1088.8334
$ws.Range("H122").Value = 4015.8276
$ws.Range("I122").Value = 2501.95
$ws.Range("J122").Value = 7380
$ws.Range("K122").Value = 7505.849999999999
$ws.Range("L122").Value = 22140
$ws.Range("M122").Value = -5055.849999999999
$ws.Range("N122").Value = -27040
$ws.Range("H132").Value = 5210.587
$ws.Range("I132").Value = 2506.074
$ws.Range("K132").Value = 7518.222
$ws.Range("M132").Value = -4988.222

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4445936
$ws.Range("I132").Value = 901.44446
$ws.Range("J132").Value = 15876025
$ws.Range("K132").Value = 2704.33338
$ws.Range("L132").Value = 47628075
$ws.Range("M132").Value = -174.33338
$ws.Range("N132").Value = -47633135
